$wb = $excel.ActiveWorkbook

# Sheet: 展览
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibitionUpdates = @{
    "F2" = 1130
    "F4" = 268
    "F5" = 1808
    "F7" = 548
    "F8" = 548
    "F9" = 5061
    "F13" = 997
    "F14" = 339
    "F15" = 1306
    "F18" = 3020
    "F19" = 1862
    "F21" = 51
    "F23" = 89
    "F24" = 651
    "F25" = 953
    "F26" = 320
    "F28" = 3341
    "F29" = 1052
    "F30" = 2593
    "F31" = 264
    "F32" = 1624
    "F33" = 3789
    "F34" = 102
    "F35" = 906
    "F36" = 440
    "F37" = 1184
    "F38" = 13
    "F39" = 953
    "F40" = 1215
    "F41" = 36
    "F42" = 928
    "F43" = 603
    "F44" = 390
    "F46" = 297
    "F47" = 3518
}
foreach ($cellRef in $wsExhibitionUpdates.Keys) {
    $wsExhibition.Range($cellRef).Value = $wsExhibitionUpdates[$cellRef]
}

# Sheet: 演出
$wsPerformance = $wb.Worksheets.Item("演出")
$wsPerformanceUpdates = @{
    "F11" = 894
}
foreach ($cellRef in $wsPerformanceUpdates.Keys) {
    $wsPerformance.Range($cellRef).Value = $wsPerformanceUpdates[$cellRef]
}

# Sheet: 全部类型
$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypesUpdates = @{
    "F2" = 1130
    "F4" = 268
    "F6" = 1808
    "F8" = 548
    "F9" = 548
    "F10" = 5061
    "F15" = 339
    "F16" = 1306
    "F17" = 3020
    "F19" = 1862
    "F21" = 51
    "F23" = 894
    "F25" = 89
    "F26" = 953
    "F27" = 320
    "F28" = 3342
    "F30" = 1052
    "F32" = 2593
    "F33" = 1624
    "F34" = 3789
    "F36" = 102
    "F37" = 906
    "F38" = 1184
    "F39" = 13
    "F40" = 953
    "F42" = 1215
    "F43" = 36
    "F44" = 928
    "F45" = 603
    "F48" = 297
    "F49" = 3518
}
foreach ($cellRef in $wsAllTypesUpdates.Keys) {
    $wsAllTypes.Range($cellRef).Value = $wsAllTypesUpdates[$cellRef]
}
